{"js": "// Replace every math-fact cell in the 20x5 practice table (row-major\n// order) with the new values from the target revision. The title\n// paragraph above the table is left untouched.\nconst newValues = [\n    \"60-37=\",\n    \"22+41=\",\n    \"16+50=\",\n    \"56-31=\",\n    \"88-68=\",\n    \"74-31=\",\n    \"45+52=\",\n    \"18+35=\",\n    \"41+49=\",\n    \"30-21=\",\n    \"0+26=\",\n    \"1+55=\",\n    \"58+31=\",\n    \"34+40=\",\n    \"31+52=\",\n    \"70-19=\",\n    \"10+75=\",\n    \"59+37=\",\n    \"38+21=\",\n    \"18+81=\",\n    \"68+5=\",\n    \"75-11=\",\n    \"53+0=\",\n    \"61-7=\",\n    \"85-21=\",\n    \"76-68=\",\n    \"55-50=\",\n    \"26+2=\",\n    \"96-44=\",\n    \"19+49=\",\n    \"6+53=\",\n    \"28+32=\",\n    \"45+28=\",\n    \"78-64=\",\n    \"80+5=\",\n    \"57+16=\",\n    \"3+12=\",\n    \"43+18=\",\n    \"50-18=\",\n    \"19-10=\",\n    \"61+27=\",\n    \"92-81=\",\n    \"93-22=\",\n    \"96-40=\",\n    \"38+59=\",\n    \"65+0=\",\n    \"58-12=\",\n    \"55-50=\",\n    \"95-39=\",\n    \"2+28=\",\n    \"78+1=\",\n    \"28+33=\",\n    \"2+10=\",\n    \"15+9=\",\n    \"67-67=\",\n    \"78-54=\",\n    \"49+50=\",\n    \"47+45=\",\n    \"84-2=\",\n    \"83-17=\",\n    \"22-0=\",\n    \"60+10=\",\n    \"28+8=\",\n    \"65+0=\",\n    \"2+26=\",\n    \"33-30=\",\n    \"79-38=\",\n    \"67-32=\",\n    \"25+26=\",\n    \"39-37=\",\n    \"34+41=\",\n    \"54-31=\",\n    \"69+12=\",\n    \"84-4=\",\n    \"1+79=\",\n    \"87-31=\",\n    \"69-66=\",\n    \"47-34=\",\n    \"7+80=\",\n    \"86-63=\",\n    \"18-14=\",\n    \"12+42=\",\n    \"83+3=\",\n    \"33+36=\",\n    \"56-18=\",\n    \"29-14=\",\n    \"22+23=\",\n    \"16+66=\",\n    \"60-26=\",\n    \"81-36=\",\n    \"95-50=\",\n    \"93-32=\",\n    \"10-2=\",\n    \"9+88=\",\n    \"35+11=\",\n    \"7+75=\",\n    \"70-67=\",\n    \"55+29=\",\n    \"6+50=\",\n    \"6+93=\"\n];\n\nconst table = context.document.body.tables.getFirst();\ntable.load(\"values\");\nawait context.sync();\n\nconst values = table.values;\nconst numRows = values.length;\nconst numCols = values[0].length;\n\nlet idx = 0;\nfor (let r = 0; r < numRows; r++) {\n    for (let c = 0; c < numCols; c++) {\n        if (idx < newValues.length) {\n            table.getCell(r, c).value = newValues[idx];\n        }\n        idx++;\n    }\n}\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n$newValues = @(\n    \"60-37=\",\n    \"22+41=\",\n    \"16+50=\",\n    \"56-31=\",\n    \"88-68=\",\n    \"74-31=\",\n    \"45+52=\",\n    \"18+35=\",\n    \"41+49=\",\n    \"30-21=\",\n    \"0+26=\",\n    \"1+55=\",\n    \"58+31=\",\n    \"34+40=\",\n    \"31+52=\",\n    \"70-19=\",\n    \"10+75=\",\n    \"59+37=\",\n    \"38+21=\",\n    \"18+81=\",\n    \"68+5=\",\n    \"75-11=\",\n    \"53+0=\",\n    \"61-7=\",\n    \"85-21=\",\n    \"76-68=\",\n    \"55-50=\",\n    \"26+2=\",\n    \"96-44=\",\n    \"19+49=\",\n    \"6+53=\",\n    \"28+32=\",\n    \"45+28=\",\n    \"78-64=\",\n    \"80+5=\",\n    \"57+16=\",\n    \"3+12=\",\n    \"43+18=\",\n    \"50-18=\",\n    \"19-10=\",\n    \"61+27=\",\n    \"92-81=\",\n    \"93-22=\",\n    \"96-40=\",\n    \"38+59=\",\n    \"65+0=\",\n    \"58-12=\",\n    \"55-50=\",\n    \"95-39=\",\n    \"2+28=\",\n    \"78+1=\",\n    \"28+33=\",\n    \"2+10=\",\n    \"15+9=\",\n    \"67-67=\",\n    \"78-54=\",\n    \"49+50=\",\n    \"47+45=\",\n    \"84-2=\",\n    \"83-17=\",\n    \"22-0=\",\n    \"60+10=\",\n    \"28+8=\",\n    \"65+0=\",\n    \"2+26=\",\n    \"33-30=\",\n    \"79-38=\",\n    \"67-32=\",\n    \"25+26=\",\n    \"39-37=\",\n    \"34+41=\",\n    \"54-31=\",\n    \"69+12=\",\n    \"84-4=\",\n    \"1+79=\",\n    \"87-31=\",\n    \"69-66=\",\n    \"47-34=\",\n    \"7+80=\",\n    \"86-63=\",\n    \"18-14=\",\n    \"12+42=\",\n    \"83+3=\",\n    \"33+36=\",\n    \"56-18=\",\n    \"29-14=\",\n    \"22+23=\",\n    \"16+66=\",\n    \"60-26=\",\n    \"81-36=\",\n    \"95-50=\",\n    \"93-32=\",\n    \"10-2=\",\n    \"9+88=\",\n    \"35+11=\",\n    \"7+75=\",\n    \"70-67=\",\n    \"55+29=\",\n    \"6+50=\",\n    \"6+93=\"\n)\n\n$table = $d.Tables.Item(1)\n$rowCount = $table.Rows.Count\n$colCount = $table.Columns.Count\n\n$idx = 0\nfor ($r = 1; $r -le $rowCount; $r++) {\n    for ($c = 1; $c -le $colCount; $c++) {\n        $cell = $table.Cell($r, $c)\n        $cell.Range.Text = $newValues[$idx]\n        $idx++\n    }\n}\n"}
